# "travel risk, sorting by time, badge"
# Adds two new entries (rows) to the daily sheet's work log:
#   - row 41: continuation of the existing day (job #3) - "Emails with Hassan and Dr."
#   - row 42: a new day (43924 = 2020-04-03, job #1) - "Travel Risk : Low/Medium/High,
#             Sort by travel time, shortest badge"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 41 - same day as row 39/40 (no date in column A), job no. "3"
$ws.Range("B41").Value = "3"
$ws.Range("C41").Value = "Emails with Hassan and Dr."
$ws.Range("D41").Value = 0.79166666666666663
$ws.Range("E41").Value = 0.8125
$ws.Range("F41").Value = 0.5

# Row 42 - new day, 43924 (2020-04-03), job no. "1"
$ws.Range("A42").Value = 43924
$ws.Range("B42").Value = "1"
$ws.Range("C42").Value = "Travel Risk : Low/Medium/High, Sort by travel time, shortest badge"
$ws.Range("D42").Value = 0.875
$ws.Range("E42").Value = 0.9375
$ws.Range("F42").Value = 1.5

# Leave the cursor/selection on the newly added last cell, scrolled so the
# new rows are in view - matches the author's saved view state.
$ws.Range("F42").Select()
